$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Update existing F/G values for rows 625-658 ---
$ws.Range("F625").Value = 43713

$ws.Range("F631").Value = 41803
$ws.Range("F632").Value = 44078

$ws.Range("F635").Value = 82997
$ws.Range("F636").Value = 49786
$ws.Range("F637").Value = 43348
$ws.Range("F638").Value = 37361
$ws.Range("F639").Value = 40414

$ws.Range("F640").Value = 19666
$ws.Range("G640").Value = 1226

$ws.Range("F641").Value = 33625
$ws.Range("G641").Value = 1370

$ws.Range("F642").Value = 67280

$ws.Range("F649").Value = 61901

$ws.Range("F652").Value = 34303
$ws.Range("F653").Value = 33128
$ws.Range("F654").Value = 13757

$ws.Range("F655").Value = 23953
$ws.Range("G655").Value = 770

$ws.Range("F656").Value = 49150
$ws.Range("G656").Value = 1171

$ws.Range("F657").Value = 32447
$ws.Range("G657").Value = 811

$ws.Range("F658").Value = 25386
$ws.Range("G658").Value = 703

# --- Append new rows 659-662 ---
$ws.Range("A659").Value = 44553
$ws.Range("A659").NumberFormat = "yyyy-mm-dd"
$ws.Range("B659").Value = 824172
$ws.Range("C659").Value = 17934
$ws.Range("D659").Value = 4679
$ws.Range("E659").Value = 16290
$ws.Range("F659").Value = 21002
$ws.Range("G659").Value = 674

$ws.Range("A660").Value = 44554
$ws.Range("A660").NumberFormat = "yyyy-mm-dd"
$ws.Range("B660").Value = 825246
$ws.Range("C660").Value = 4201
$ws.Range("D660").Value = 1074
$ws.Range("E660").Value = 16336
$ws.Range("F660").Value = 4982
$ws.Range("G660").Value = 181

$ws.Range("A661").Value = 44555
$ws.Range("A661").NumberFormat = "yyyy-mm-dd"
$ws.Range("B661").Value = 825629
$ws.Range("C661").Value = 2376
$ws.Range("D661").Value = 383
$ws.Range("E661").Value = 16398
$ws.Range("F661").Value = 3614
$ws.Range("G661").Value = 245

$ws.Range("A662").Value = 44556
$ws.Range("A662").NumberFormat = "yyyy-mm-dd"
$ws.Range("B662").Value = 826850
$ws.Range("C662").Value = 4995
$ws.Range("D662").Value = 1221
$ws.Range("E662").Value = 16445
$ws.Range("F662").Value = 7455
$ws.Range("G662").Value = 358
